$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "30.318.09"
$ws.Range("E2").Value = "  -0.39%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.935.45"
$ws.Range("E3").Value = "  -0.04%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "'0.9989"
$ws.Range("E4").Value = "  -0.82%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'250.80"
$ws.Range("E5").Value = "  +1.25%  "

# Row 6 - XRP
$ws.Range("D6").Value = "'0.7265"
$ws.Range("E6").Value = "  +4.64%  "

# Row 7 - USDC
$ws.Range("D7").Value = "'0.9988"
$ws.Range("E7").Value = "  -0.86%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.3312"
$ws.Range("E8").Value = "  +2.24%  "

# Row 9 - Solana
$ws.Range("D9").Value = "'28.01"
$ws.Range("E9").Value = "  +5.54%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "'0.07295"
$ws.Range("E10").Value = "  +7.51%  "

# Row 11 - Polygon
$ws.Range("D11").Value = "'0.8080"
$ws.Range("E11").Value = "  +1.64%  "

# Row 12 - TRON (only price changes)
$ws.Range("D12").Value = "'0.08098"

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.935.10"
$ws.Range("E13").Value = "  -0.11%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'5.480"
$ws.Range("E14").Value = "  +1.69%  "

# Row 15 - Litecoin
$ws.Range("D15").Value = "'94.82"
$ws.Range("E15").Value = "  +0.71%  "

# Row 16 - Avalanche
$ws.Range("D16").Value = "'15.13"
$ws.Range("E16").Value = "  +4.48%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "30.303.53"
$ws.Range("E17").Value = "  -0.47%  "

# Row 18 - was ShibaInu, now BitcoinCash
$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").Value = "'253.55"
$ws.Range("E18").Value = "  -2.93%  "

# Row 19 - was BitcoinCash, now ShibaInu
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000008225"
$ws.Range("E19").Value = "  +5.35%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "'5.814"
$ws.Range("E20").Value = "  -0.94%  "

# Row 21 - WrappedliquidstakedEther2.0
$ws.Range("D21").Value = "2.186.11"
$ws.Range("E21").Value = "  -0.78%  "

# Row 22 - Dai
$ws.Range("D22").Value = "'0.9986"
$ws.Range("E22").Value = "  -0.87%  "

# Row 23
$ws.Range("D23").Value = "'0.9988"
$ws.Range("E23").Value = "  -0.78%  "

# Row 24
$ws.Range("D24").Value = "'6.967"
$ws.Range("E24").Value = "  +2.02%  "

# Row 25
$ws.Range("D25").Value = "'9.774"
$ws.Range("E25").Value = "  +1.71%  "

# Row 26
$ws.Range("D26").Value = "'165.35"
$ws.Range("E26").Value = "  +4.21%  "

# Row 27
$ws.Range("D27").Value = "'2.360"
$ws.Range("E27").Value = "  +4.80%  "

# Row 28
$ws.Range("D28").Value = "'19.34"
$ws.Range("E28").Value = "  +2.78%  "

# Row 29
$ws.Range("D29").Value = "'0.1309"
$ws.Range("E29").Value = "  +0.20%  "

# Row 30
$ws.Range("D30").Value = "'1.353"
$ws.Range("E30").Value = "  -1.02%  "

# Row 31
$ws.Range("D31").Value = "'1.540"
$ws.Range("E31").Value = "  -1.25%  "

# Row 32
$ws.Range("D32").Value = "'4.438"
$ws.Range("E32").Value = "  +0.57%  "

# Row 33
$ws.Range("D33").Value = "'4.201"
$ws.Range("E33").Value = "  -0.55%  "

# Row 34
$ws.Range("D34").Value = "'0.05250"
$ws.Range("E34").Value = "  +3.23%  "

# Row 35
$ws.Range("D35").Value = "'1.273"
$ws.Range("E35").Value = "  +6.59%  "

# Row 36
$ws.Range("D36").Value = "'0.7505"
$ws.Range("E36").Value = "  +0.23%  "

# Row 37
$ws.Range("D37").Value = "'2.759"
$ws.Range("E37").Value = "  +1.10%  "

# Row 38 - VeChain (only volume changes)
$ws.Range("E38").Value = "  +2.54%  "

# Row 39
$ws.Range("D39").Value = "'2.810"
$ws.Range("E39").Value = "  +1.02%  "

# Row 40
$ws.Range("D40").Value = "'79.25"
$ws.Range("E40").Value = "  -0.89%  "

# Row 41
$ws.Range("D41").Value = "'6.443"
$ws.Range("E41").Value = "  -1.76%  "

# Row 42
$ws.Range("D42").Value = "'0.4546"
$ws.Range("E42").Value = "  +2.72%  "

# Row 43
$ws.Range("D43").Value = "'2.038"
$ws.Range("E43").Value = "  -0.59%  "

# Row 44
$ws.Range("D44").Value = "'0.8447"
$ws.Range("E44").Value = "  +0.27%  "

# Row 45
$ws.Range("D45").Value = "'0.9990"
$ws.Range("E45").Value = "  -0.98%  "

# Row 46
$ws.Range("D46").Value = "'101.78"
$ws.Range("E46").Value = "  -0.13%  "

# Row 47
$ws.Range("D47").Value = "'9.752"
$ws.Range("E47").Value = "  +0.43%  "

# Row 48
$ws.Range("D48").Value = "'7.458"
$ws.Range("E48").Value = "  +2.19%  "

# Row 49
$ws.Range("D49").Value = "'36.82"
$ws.Range("E49").Value = "  +2.12%  "

# Row 50
$ws.Range("D50").Value = "'0.4195"
$ws.Range("E50").Value = "  +2.70%  "

# Row 51
$ws.Range("D51").Value = "'0.06041"
$ws.Range("E51").Value = "  +1.61%  "
